# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" between "2021-Q4" and "总计", with the
#    per-fund holdings table for the new quarter.
# 2. Update the "总计" (summary) sheet: insert a new row for "2022-Q1" above
#    the existing "2021-Q4" row.
#
# Implementation note: rather than inserting a brand-new blank sheet (which
# would come out with Excel's default page setup / sheetFormatPr instead of
# matching the rest of the workbook), we duplicate the existing "总计" sheet
# first. One copy is repurposed (cleared + repopulated) as "2022-Q1"; the
# other keeps being "总计" and just gets the extra summary row inserted. That
# way both end up sharing the original sheet's formatting, same as the
# source workbook.

$wb = $excel.ActiveWorkbook
$wsQ4 = $wb.Worksheets.Item(1)   # "2021-Q4" -- position never changes; used
                                 # purely as a source of formatting to copy.

$wsSumOrig = $wb.Worksheets.Item(2)   # "总计", about to be duplicated

# Duplicate "总计" right after itself -> position 2 keeps the original sheet
# (same underlying file/formatting), position 3 is the new copy.
$wsSumOrig.Copy($null, $wsSumOrig)

# NOTE: Worksheets.Item(N) is position-based. From here on, re-fetch sheets
# by position fresh instead of reusing references obtained before the Copy.
$wsQ1 = $wb.Worksheets.Item(2)
$wsQ1.Name = "2022-Q1"

$wsSum = $wb.Worksheets.Item(3)
$wsSum.Name = "总计"

# ---------------------------------------------------------------------
# 1. Rebuild the "2022-Q1" sheet's content (currently still holds the old
#    "总计" header/row, since it was cloned before being renamed)
# ---------------------------------------------------------------------
$wsQ1.Cells.Clear()

$wsQ1.Range("B1").Value = "基金代码"
$wsQ1.Range("C1").Value = "基金名称"
$wsQ1.Range("D1").Value = "基金规模"
$wsQ1.Range("E1").Value = "股票总仓位"
$wsQ1.Range("F1").Value = "仓位占比"
$wsQ1.Range("G1").Value = "持有市值(亿元)"
$wsQ1.Range("H1").Value = "仓位排名"

# Copy the header styling (bold + border + centered) from the "2021-Q4" sheet
$wsQ4.Range("B1:H1").Copy()
$wsQ1.Range("B1:H1").PasteSpecial(-4122)

# Data rows: A (index, numeric, styled) / B..G (text, even when numeric-looking,
# to preserve leading zeros / match source formatting) / H (plain number)
$rows = @(
    @(0, "009230", "鹏华安和混合A",         "14.02", "34.45", "1.14", "0.1598", 10),
    @(1, "009667", "鹏华安庆混合A",         "11.22", "38.92", "1.30", "0.1459", 8),
    @(2, "009231", "鹏华安和混合C",         "5.33",  "34.45", "1.14", "0.0608", 10),
    @(3, "003165", "鹏华弘嘉灵活配置混合A", "1.53",  "93.95", "3.43", "0.0525", 4),
    @(4, "004558", "汇安丰裕灵活配置混合A", "0.99",  "83.41", "3.90", "0.0386", 5),
    @(5, "009668", "鹏华安庆混合C",         "2.36",  "38.92", "1.30", "0.0307", 8),
    @(6, "003166", "鹏华弘嘉灵活配置混合C", "0.56",  "93.95", "3.43", "0.0192", 4),
    @(7, "004559", "汇安丰裕灵活配置混合C", "0.01",  "83.41", "3.90", "0.0004", 5)
)

$r = 2
foreach ($row in $rows) {
    $wsQ1.Cells.Item($r, 1).Value = $row[0]

    $wsQ1.Range("B$r`:G$r").NumberFormat = "@"
    $wsQ1.Cells.Item($r, 2).Value = $row[1]
    $wsQ1.Cells.Item($r, 3).Value = $row[2]
    $wsQ1.Cells.Item($r, 4).Value = $row[3]
    $wsQ1.Cells.Item($r, 5).Value = $row[4]
    $wsQ1.Cells.Item($r, 6).Value = $row[5]
    $wsQ1.Cells.Item($r, 7).Value = $row[6]
    $wsQ1.Range("B$r`:G$r").Style = "Normal"

    $wsQ1.Cells.Item($r, 8).Value = $row[7]

    $r = $r + 1
}

# Style column A (the index column) the same way as on the "2021-Q4" sheet
$wsQ4.Range("A2").Copy()
$wsQ1.Range("A2:A9").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. "总计" sheet: insert the new "2022-Q1" row above "2021-Q4"
# ---------------------------------------------------------------------
$wsSum.Rows.Item(2).Insert()

# Re-apply the index-column style to the new A2 (Insert leaves it unstyled)
$wsQ4.Range("A2").Copy()
$wsSum.Range("A2").PasteSpecial(-4122)

# Clear the stray formatting Insert left on B2:D2 by borrowing the (unstyled)
# format from the row below, which still holds the original "2021-Q4" data
$wsSum.Range("B3:D3").Copy()
$wsSum.Range("B2:D2").PasteSpecial(-4122)

$wsSum.Cells.Item(2, 1).Value = 0
$wsSum.Cells.Item(2, 2).Value = "2022-Q1"
$wsSum.Cells.Item(2, 3).Value = 8
$wsSum.Cells.Item(2, 4).Value = 0.51

# Bump the index of the (now shifted-down) "2021-Q4" row
$wsSum.Cells.Item(3, 1).Value = 1
